# Correct reference to Karsholt & Razowski
# Add the Karsholt & Razowski reference (referenceID 979) to the
# distribution rows for the European countries this reference covers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("region")

$rows = @(16, 23, 58, 76, 83, 107, 110, 129, 130, 138, 156, 179, 184, 204, 219, 238)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 979
}

# Leave the view where Excel would land after entering the last value:
# scrolled down so row 211 is at the top, with the cell below the last
# data row selected.
$excel.ActiveWindow.ScrollRow = 211
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C252").Select()
